$wb = $excel.ActiveWorkbook

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 33336266
$ws.Range("I113").Value = 50002052
$ws.Range("J113").Value = 4697.6
$ws.Range("K113").Value = 50002052
$ws.Range("L113").Value = 4697.6
$ws.Range("M113").Value = -49998798
$ws.Range("N113").Value = -11205.6

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3206591
$ws.Range("I141").Value = 1448.2
$ws.Range("J141").Value = 6174316
$ws.Range("K141").Value = 4344.6
$ws.Range("L141").Value = 18522948
$ws.Range("M141").Value = 835.3999999999996
$ws.Range("N141").Value = -18533308

# ARM row 44
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 31937.666
$ws.Range("J44").Value = 31937.666
$ws.Range("L44").Value = 31937.666
$ws.Range("N44").Value = -32913.666

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2913.6
$ws.Range("I45").Value = 2350.1428
$ws.Range("J45").Value = 3630.7273
$ws.Range("K45").Value = 2350.1428
$ws.Range("L45").Value = 3630.7273
$ws.Range("M45").Value = -1973.1428
$ws.Range("N45").Value = -4384.7273

# ARM row 55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 19614.5
$ws.Range("J55").Value = 19614.5
$ws.Range("L55").Value = 19614.5
$ws.Range("N55").Value = -20244.5

# ARM row 80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 31267.143
$ws.Range("J80").Value = 34845
$ws.Range("L80").Value = 34845
$ws.Range("N80").Value = -36841

# ARM row 83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 31267.143
$ws.Range("J83").Value = 34845
$ws.Range("L83").Value = 104535
$ws.Range("N83").Value = -114519

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 142859410
$ws.Range("I102").Value = 2333.3333
$ws.Range("J102").Value = 250002200
$ws.Range("K102").Value = 2333.3333
$ws.Range("L102").Value = 250002200
$ws.Range("M102").Value = -711.3332999999998
$ws.Range("N102").Value = -250005444

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1642.25
$ws.Range("I110").Value = 1800.0869
$ws.Range("J110").Value = 1238.8889
$ws.Range("K110").Value = 1800.0869
$ws.Range("L110").Value = 1238.8889
$ws.Range("M110").Value = 244.9131
$ws.Range("N110").Value = -5328.8889

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2036.2
$ws.Range("I122").Value = 1818
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 5454
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -3004
$ws.Range("N122").Value = -16900

# ARM row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 665770
$ws.Range("J123").Value = 665770
$ws.Range("L123").Value = 665770
$ws.Range("N123").Value = -675570

# BSM row 19
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 8000
$ws.Range("J19").Value = 8000
$ws.Range("L19").Value = 8000
$ws.Range("N19").Value = -8346

# BSM row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 25008.666
$ws.Range("J35").Value = 25008.666
$ws.Range("L35").Value = 25008.666
$ws.Range("N35").Value = -25628.666

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 16920.941
$ws.Range("J82").Value = 31769.5
$ws.Range("L82").Value = 31769.5
$ws.Range("N82").Value = -32535.5

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 16920.941
$ws.Range("J85").Value = 31769.5
$ws.Range("L85").Value = 31769.5
$ws.Range("N85").Value = -34421.5

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3914.2856
$ws.Range("I99").Value = 1633.3334
$ws.Range("J99").Value = 5625
$ws.Range("K99").Value = 1633.3334
$ws.Range("L99").Value = 5625
$ws.Range("M99").Value = -135.3334
$ws.Range("N99").Value = -8621

# BSM row 103
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 43943.5
$ws.Range("J103").Value = 43943.5
$ws.Range("L103").Value = 43943.5
$ws.Range("N103").Value = -46287.5

# BSM row 122
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 51980
$ws.Range("J122").Value = 51980
$ws.Range("L122").Value = 51980
$ws.Range("N122").Value = -61780

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19283.8
$ws.Range("J41").Value = 22615
$ws.Range("L41").Value = 22615
$ws.Range("N41").Value = -23471

# CRP row 55
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 5936.5
$ws.Range("I55").Value = 5873
$ws.Range("K55").Value = 5873
$ws.Range("M55").Value = -5558

# CRP row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 20350.75
$ws.Range("J60").Value = 25134.334
$ws.Range("L60").Value = 25134.334
$ws.Range("N60").Value = -26156.334

# CRP row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 27965
$ws.Range("J68").Value = 27965
$ws.Range("L68").Value = 27965
$ws.Range("N68").Value = -29463

# CRP row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 27965
$ws.Range("J71").Value = 27965
$ws.Range("L71").Value = 83895
$ws.Range("N71").Value = -91383

# CRP row 100
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 37758.57
$ws.Range("J100").Value = 37758.57
$ws.Range("L100").Value = 37758.57
$ws.Range("N100").Value = -39922.57

# CRP row 125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 54980
$ws.Range("J125").Value = 54980
$ws.Range("L125").Value = 54980
$ws.Range("N125").Value = -59900

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2117.196
$ws.Range("J68").Value = 2143.3794
$ws.Range("L68").Value = 6430.138199999999
$ws.Range("N68").Value = -8052.138199999999

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2117.196
$ws.Range("J71").Value = 2143.3794
$ws.Range("L71").Value = 19290.4146
$ws.Range("N71").Value = -27402.4146

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 901.22644
$ws.Range("I107").Value = 591.4286
$ws.Range("J107").Value = 1503.6111
$ws.Range("K107").Value = 1774.2858
$ws.Range("L107").Value = 4510.8333
$ws.Range("M107").Value = 145.7142000000001
$ws.Range("N107").Value = -8350.8333

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 997.9783
$ws.Range("J131").Value = 1009.04443
$ws.Range("L131").Value = 3027.13329
$ws.Range("N131").Value = -13107.13329

# GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 9490
$ws.Range("J5").Value = 9490
$ws.Range("L5").Value = 9490
$ws.Range("N5").Value = -9714

# GSM row 64
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 37592
$ws.Range("J64").Value = 37592
$ws.Range("L64").Value = 37592
$ws.Range("N64").Value = -38088

# GSM row 67
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H67").Value = 37592
$ws.Range("J67").Value = 37592
$ws.Range("L67").Value = 37592
$ws.Range("N67").Value = -39308

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6401.3
$ws.Range("I122").Value = 4876.625
$ws.Range("K122").Value = 14629.875
$ws.Range("M122").Value = -12179.875

# LTW row 39
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 12000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 12000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 12000
$ws.Range("M39").Value = $null
$ws.Range("N39").Value = -12920

# LTW row 62
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 34062.25
$ws.Range("J62").Value = 34062.25
$ws.Range("L62").Value = 34062.25
$ws.Range("N62").Value = -35310.25

# LTW row 65
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H65").Value = 34062.25
$ws.Range("J65").Value = 34062.25
$ws.Range("L65").Value = 102186.75
$ws.Range("N65").Value = -108426.75

# LTW row 111
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 41966.668
$ws.Range("J111").Value = 41966.668
$ws.Range("L111").Value = 41966.668
$ws.Range("N111").Value = -50146.668

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3003.2222
$ws.Range("I122").Value = 2084.75
$ws.Range("J122").Value = 3738
$ws.Range("K122").Value = 6254.25
$ws.Range("L122").Value = 11214
$ws.Range("M122").Value = -3804.25
$ws.Range("N122").Value = -16114

# WVR row 27
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 37500
$ws.Range("J27").Value = 37500
$ws.Range("L27").Value = 37500
$ws.Range("N27").Value = -37638

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 7661.6
$ws.Range("J54").Value = 7661.6
$ws.Range("L54").Value = 7661.6
$ws.Range("N54").Value = -8701.6

# WVR row 115
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 48000
$ws.Range("J115").Value = 48000
$ws.Range("L115").Value = 48000
$ws.Range("N115").Value = -51134

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3842.6316
$ws.Range("I122").Value = 3791.8333
$ws.Range("J122").Value = 3929.7144
$ws.Range("K122").Value = 11375.4999
$ws.Range("L122").Value = 11789.1432
$ws.Range("M122").Value = -8925.499899999999
$ws.Range("N122").Value = -16689.1432
